# Update labels in the "Graduate Students" sheet: remove duplicated text
# (e.g. "MaleMale" -> "Male") and strip the leading indentation whitespace
# that was baked into the text. Row 2 is consolidated into "All students".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Graduate Students")

$labels = @{
    2  = "All students"
    3  = "Male"
    4  = "Female"
    5  = "U.S. citizens and permanent residents"
    6  = "Hispanic or Latino"
    7  = "Not Hispanic or Latino"
    8  = "American Indian or Alaska Native"
    9  = "Asian"
    10 = "Black or African American"
    11 = "Native Hawaiian or Other Pacific Islander"
    12 = "White"
    13 = "More than one race"
    14 = "Unknown ethnicity and race"
    15 = "Temporary visa holders"
    16 = "Science and engineering"
    17 = "Science"
    18 = "Agricultural and veterinary sciences"
    19 = "Biological and biomedical sciences"
    20 = "Communication"
    21 = "Computer and information sciences"
    22 = "Family and consumer sciences and human sciences"
    23 = "Geosciences, atmospheric sciences, and ocean sciences"
    24 = "Mathematics and statistics"
    25 = "Multidisciplinary and interdisciplinary studies"
    26 = "Natural resources and conservation"
    27 = "Neurobiology and neuroscience"
    28 = "Physical sciences"
    29 = "Psychology"
    30 = "Social sciences"
    31 = "Engineering"
    32 = "Aerospace, aeronautical, and astronautical engineering"
    33 = "Biological, biomedical, and biosystems engineering"
    34 = "Chemical, petroleum, and chemical-related engineering"
    35 = "Civil, environmental, transportation and related engineering fields"
    36 = "Electrical, electronics, communications and computer engineering"
    37 = "Industrial, manufacturing, systems engineering and operations research"
    38 = "Mechanical engineering"
    39 = "Metallurgical, mining, materials and related engineering fields"
    40 = "Other engineering"
    41 = "Health"
    42 = "Clinical medicine"
    43 = "Other health"
}

foreach ($row in $labels.Keys) {
    $ws.Cells.Item($row, 1).Value = $labels[$row]
}
